# Updated cryptos list on Mon Nov 20 07:29:48 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.174.02"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.002.48"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.57"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.58"
$ws.Range("E7").Value = "  -1.88%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.01"
$ws.Range("E12").Value = "  +5.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.27"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.298.98"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.843"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.001.78"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.055.11"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.17"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.94"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.33"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.62"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  +12.97%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.79"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +5.57%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +6.22%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.31"
$ws.Range("E39").Value = "  -5.21%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0213"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.55"
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.99"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.367.40"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.37"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +13.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.85"
$ws.Range("E50").Value = "  +5.53%  "
$ws.Range("E51").Value = "  +0.04%  "
